$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (peak_Normalized_Measurement_p)
$ws.Range("B5").Value = 0.114
$ws.Range("C5").Value = 0.012
$ws.Range("D5").Value = 0

# Row 6 (peak_Normalized_Measurement_ci)
$ws.Range("B6").Value = "(-0.078, 1.065)"
$ws.Range("C6").Value = "(1.493, 3.942)"
$ws.Range("D6").Value = "(10.816, 54.506)"

# Row 11 (peak_Change_Rate_p)
$ws.Range("C11").Value = 0.003
$ws.Range("D11").Value = 0.01

# Row 12 (peak_Change_Rate_ci)
$ws.Range("B12").Value = "(-0.649, 0.369)"
$ws.Range("C12").Value = "(-0.652, 0.246)"
$ws.Range("D12").Value = "(1.023, 14.155)"
